# Se agrega funcionalidad para limpiar los movimientos de stock
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 41 ("Magic Button para borrar movimientos de stock") is now done:
# replace the "en proceso" text with a 100% (percentage-styled) value,
# matching the "done" convention used elsewhere on the sheet.
$ws.Range("C41").Value = 1
$ws.Range("C41").NumberFormat = "0%"

# New task in row 44 ("en stock agregar codigo - descripcion"): assign it
# to Agustina and mark it "en proceso".
$ws.Range("B44").Value = "Agustina"
$ws.Range("C44").Value = "en proceso"

# Move the active selection down to C45.
$ws.Range("C45").Select()
